# Apply the "Custody Status Change" IEPD documentation edit:
#  - Insert a new row for "Person Ethnicity Code (Pima County)" right after the
#    existing "Ethnicity" row (new row 21), mirroring the pattern used by the
#    neighbouring "Race Code (Pima County)" row.
#  - Insert a new row for "Charge Disposition" right after "Charge Description"
#    (new row 75, pushing "Statute or Ordinance Number" etc. down).
#  - Update the saved window position and the active-cell selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New row 21: Person Ethnicity Code (Pima County) ------------------
$ws.Rows("21").Insert()

$ws.Range("A21").Value = "x-ext"
$ws.Range("B21").Value = "Ethnicity"
$ws.Range("C21").Value = "Person Ethnicity Code (Pima County)"
$ws.Range("D21").Value = ""
$ws.Range("E21").Value = "/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/nc:Person[@structures:id=/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/j:Booking/j:BookingSubject/nc:RoleOfPerson/@structures:ref]/pc-bkg-codes:PersonEthnicityCode"

$rng21 = $ws.Range("A21:E21")
$rng21.HorizontalAlignment = -4131
$rng21.VerticalAlignment = -4160
$rng21.WrapText = $true
$ws.Rows("21").RowHeight = 30

# --- 2. New row 75: Charge Disposition ------------------------------------
$ws.Rows("75").Insert()

$ws.Range("B75").Value = "Charge Disposition"
$ws.Range("C75").Value = "ChargeDisposition Text"
$ws.Range("E75").Value = "/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/j:Charge[@structures:id=/cscr-doc:CustodyStatusChangeReport/cscr-ext:Custody/j:Arrest/j:ArrestCharge/@structures:ref]/j:ChargeDisposition/nc:DispositionText"

$rng75 = $ws.Range("A75:E75")
$rng75.HorizontalAlignment = -4131
$rng75.VerticalAlignment = -4160
$rng75.WrapText = $true
$ws.Rows("75").RowHeight = 30

# --- 3. Selection / window bookkeeping ------------------------------------
$ws.Range("E75").Select()

$wb.Windows.Item(1).WindowState = $wb.Windows.Item(1).WindowState
